# Update "paises" workbook: refresh COVID figures for a set of countries
# and bump the "Datos actualizados..." timestamp, per the commit
# "Update countries & provincias Spain".
#
# The per-country row order in the "Pais" sheet does not change; only the
# numeric Casos totales / Nuevos casos / Casos activos / Recuperados /
# Casos criticos / Muertes hoy / Muertes columns (B:H) for the affected
# rows are refreshed, plus the timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Pais")

# Timestamp update (A1): 16:22 -> 16:52
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 16:52"

# Row -> new [B, C, D, E, F, G, H] values (Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$updates = @{
    4   = @(402471, 2136, 21804, 367753, 9195, 73, 12914)
    27  = @(5546, 430, 1115, 4383, 362, 5, 48)
    28  = @(5480, 129, 468, 4848, 0, 4, 164)
    29  = @(5386, 315, 1621, 3547, 127, 15, 218)
    34  = @(4194, 159, 467, 3667, 25, 3, 60)
    41  = @(2932, 137, 631, 2260, 41, 0, 41)
    49  = @(2111, 155, 50, 1953, 147, 10, 108)
    62  = @(1202, 80, 452, 681, 0, 4, 69)
    63  = @(1185, 36, 72, 1089, 11, 3, 24)
    72  = @(822, 105, 63, 751, 23, 0, 8)
    73  = @(821, 10, 467, 349, 3, 0, 5)
    74  = @(803, 39, 79, 690, 4, 1, 34)
    81  = @(593, 16, 42, 527, 27, 1, 24)
    103 = @(273, 5, 19, 247, 3, 0, 7)
    104 = @(270, 42, 33, 233, 5, 0, 4)
    168 = @(17, 7, 1, 16, 0, 0, 0)
    169 = @(17, 0, 2, 13, 0, 0, 2)
    171 = @(16, 0, 3, 13, 0, 0, 0)
    172 = @(16, 1, 4, 12, 0, 0, 0)
    173 = @(15, 0, 0, 15, 0, 0, 0)
    174 = @(15, 1, 0, 15, 0, 0, 0)
    175 = @(15, 0, 1, 14, 0, 0, 0)
    176 = @(14, 0, 1, 13, 0, 0, 0)
    177 = @(14, 0, 2, 10, 0, 0, 2)
    178 = @(14, 0, 3, 8, 0, 0, 3)
    179 = @(13, 0, 7, 5, 0, 0, 1)
    180 = @(12, 0, 0, 12, 2, 0, 0)
    182 = @(11, 0, 0, 11, 0, 0, 0)
    183 = @(11, 0, 0, 9, 0, 0, 2)
    184 = @(11, 0, 10, 1, 0, 0, 0)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i   # column B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
